$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells touched below to remain plain text,
# matching the original inlineStr string storage (avoids Excel
# auto-converting numeric-looking strings like "0.556" into numbers).
$priceCells = "D2","D3","D6","D8","D9","D11","D12","D13","D14","D15","D16","D18","D20","D24","D25","D26","D30","D31","D32","D35","D36","D37","D40","D42","D44","D45","D46","D48","D49","D51"
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "34.545.67"
$ws.Range("E2").Value = "  +1.00%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.798.52"
$ws.Range("E3").Value = "  +0.82%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  +0.50%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.556"
$ws.Range("E6").Value = "  +1.43%  "

# Row 8 - Solana
$ws.Range("D8").Value = "32.69"
$ws.Range("E8").Value = "  +2.42%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.296"
$ws.Range("E9").Value = "  +1.61%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.45%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0950"
$ws.Range("E11").Value = "  +0.36%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.058.08"
$ws.Range("E12").Value = "  +0.79%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "11.19"
$ws.Range("E13").Value = "  +1.18%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.794.07"
$ws.Range("E14").Value = "  +0.32%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.639"
$ws.Range("E15").Value = "  +2.24%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "34.544.07"
$ws.Range("E16").Value = "  +1.07%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  +3.02%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "68.87"
$ws.Range("E18").Value = "  +1.31%  "

# Row 19
$ws.Range("E19").Value = "  -0.07%  "

# Row 20
$ws.Range("D20").Value = "247.17"
$ws.Range("E20").Value = "  +0.08%  "

# Row 21
$ws.Range("E21").Value = "  +2.56%  "

# Row 22
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("E23").Value = "  +2.21%  "

# Row 24
$ws.Range("D24").Value = "169.02"
$ws.Range("E24").Value = "  +3.91%  "

# Row 25
$ws.Range("D25").Value = "2.06"
$ws.Range("E25").Value = "  +0.61%  "

# Row 26
$ws.Range("D26").Value = "7.28"
$ws.Range("E26").Value = "  +1.07%  "

# Row 27
$ws.Range("E27").Value = "  +1.69%  "

# Row 28
$ws.Range("E28").Value = "  +1.92%  "

# Row 29
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("D30").Value = "4.13"
$ws.Range("E30").Value = "  +9.95%  "

# Row 31 - becomes Hedera (was Filecoin)
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.0526"
$ws.Range("E31").Value = "  +0.93%  "

# Row 32 - becomes Filecoin (was Hedera)
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.82"
$ws.Range("E32").Value = "  +2.12%  "

# Row 33
$ws.Range("E33").Value = "  +0.76%  "

# Row 34
$ws.Range("E34").Value = "  +2.47%  "

# Row 35
$ws.Range("D35").Value = "1.430.87"
$ws.Range("E35").Value = "  -0.90%  "

# Row 36
$ws.Range("D36").Value = "2.60"
$ws.Range("E36").Value = "  +8.07%  "

# Row 37
$ws.Range("D37").Value = "0.674"
$ws.Range("E37").Value = "  +2.66%  "

# Row 38
$ws.Range("E38").Value = "  +2.95%  "

# Row 39
$ws.Range("E39").Value = "  +0.19%  "

# Row 40
$ws.Range("D40").Value = "85.04"
$ws.Range("E40").Value = "  +5.88%  "

# Row 41
$ws.Range("E41").Value = "  +2.00%  "

# Row 42
$ws.Range("D42").Value = "0.938"
$ws.Range("E42").Value = "  +1.37%  "

# Row 43
$ws.Range("E43").Value = "  +3.43%  "

# Row 44
$ws.Range("D44").Value = "13.87"
$ws.Range("E44").Value = "  +2.57%  "

# Row 45
$ws.Range("D45").Value = "0.0524"
$ws.Range("E45").Value = "  +2.64%  "

# Row 46
$ws.Range("D46").Value = "6.09"
$ws.Range("E46").Value = "  +0.16%  "

# Row 47
$ws.Range("E47").Value = "  +0.84%  "

# Row 48
$ws.Range("D48").Value = "1.957.07"
$ws.Range("E48").Value = "  +0.68%  "

# Row 49
$ws.Range("D49").Value = "105.81"
$ws.Range("E49").Value = "  +1.14%  "

# Row 50
$ws.Range("E50").Value = "  -0.02%  "

# Row 51
$ws.Range("D51").Value = "0.0₆0128"
$ws.Range("E51").Value = "  -4.45%  "
